$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Problem 1
$ws.Range("B4").Value = 2500
$ws.Range("C4").Value = 95
$ws.Range("D4").Value = 2417
$ws.Range("E4").Value = 2473
$ws.Range("F4").Value = 2768
$ws.Range("G4").Value = 52
$ws.Range("H4").Value = 24
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 123
$ws.Range("L4").Value = 9193
$ws.Range("M4").Value = 830
$ws.Range("N4").Value = 7834
$ws.Range("O4").Value = 9364
$ws.Range("P4").Value = 10590

# Problem 2
$ws.Range("B5").Value = 2323
$ws.Range("C5").Value = 203
$ws.Range("D5").Value = 2148
$ws.Range("E5").Value = 2238
$ws.Range("F5").Value = 2770
$ws.Range("G5").Value = 53
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 52
$ws.Range("J5").Value = 53
$ws.Range("K5").Value = 59
$ws.Range("L5").Value = 12608
$ws.Range("M5").Value = 1109
$ws.Range("N5").Value = 10505
$ws.Range("O5").Value = 12689
$ws.Range("P5").Value = 14151

# Problem 3
$ws.Range("B6").Value = 44
$ws.Range("C6").Value = 27
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 62
$ws.Range("F6").Value = 65
$ws.Range("G6").Value = 168
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 143
$ws.Range("J6").Value = 151
$ws.Range("K6").Value = 334
$ws.Range("L6").Value = 4446
$ws.Range("M6").Value = 936
$ws.Range("N6").Value = 3624
$ws.Range("O6").Value = 4132
$ws.Range("P6").Value = 6870

# Problem 4
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 59
$ws.Range("G7").Value = 422
$ws.Range("H7").Value = 52
$ws.Range("I7").Value = 403
$ws.Range("J7").Value = 405
$ws.Range("K7").Value = 577
$ws.Range("L7").Value = 6194
$ws.Range("M7").Value = 301
$ws.Range("N7").Value = 5785
$ws.Range("O7").Value = 6290
$ws.Range("P7").Value = 6707

# Problem 5
$ws.Range("B8").Value = 1927
$ws.Range("C8").Value = 706
$ws.Range("D8").Value = 734
$ws.Range("E8").Value = 2303
$ws.Range("F8").Value = 2554
$ws.Range("G8").Value = 661
$ws.Range("H8").Value = 103
$ws.Range("I8").Value = 590
$ws.Range("J8").Value = 615
$ws.Range("K8").Value = 870
$ws.Range("L8").Value = 17088
$ws.Range("M8").Value = 1828
$ws.Range("N8").Value = 13761
$ws.Range("O8").Value = 17072
$ws.Range("P8").Value = 19766

# Problem 6
$ws.Range("B9").Value = 3197
$ws.Range("C9").Value = 581
$ws.Range("D9").Value = 2313
$ws.Range("E9").Value = 3291
$ws.Range("F9").Value = 4308
$ws.Range("G9").Value = 746
$ws.Range("H9").Value = 7
$ws.Range("I9").Value = 727
$ws.Range("J9").Value = 749
$ws.Range("K9").Value = 752
$ws.Range("L9").Value = 16826
$ws.Range("M9").Value = 2040
$ws.Range("N9").Value = 13654
$ws.Range("O9").Value = 16745
$ws.Range("P9").Value = 20358

# Problem 7
$ws.Range("B10").Value = 1005
$ws.Range("C10").Value = 677
$ws.Range("D10").Value = 212
$ws.Range("E10").Value = 747
$ws.Range("F10").Value = 2245
$ws.Range("G10").Value = 767
$ws.Range("H10").Value = 11
$ws.Range("I10").Value = 760
$ws.Range("J10").Value = 763
$ws.Range("K10").Value = 799
$ws.Range("L10").Value = 19017
$ws.Range("M10").Value = 2743
$ws.Range("N10").Value = 15291
$ws.Range("O10").Value = 18737
$ws.Range("P10").Value = 24726

# Problem 8
$ws.Range("B11").Value = 342
$ws.Range("C11").Value = 694
$ws.Range("D11").Value = 77
$ws.Range("E11").Value = 90
$ws.Range("F11").Value = 2416
$ws.Range("G11").Value = 1008
$ws.Range("H11").Value = 45
$ws.Range("I11").Value = 957
$ws.Range("J11").Value = 996
$ws.Range("K11").Value = 1095
$ws.Range("L11").Value = 9071
$ws.Range("M11").Value = 4526
$ws.Range("N11").Value = 5031
$ws.Range("O11").Value = 7366
$ws.Range("P11").Value = 18554

# Problem 9
$ws.Range("B12").Value = 574
$ws.Range("C12").Value = 900
$ws.Range("D12").Value = 40
$ws.Range("E12").Value = 118
$ws.Range("F12").Value = 2686
$ws.Range("G12").Value = 1237
$ws.Range("H12").Value = 117
$ws.Range("I12").Value = 1088
$ws.Range("J12").Value = 1313
$ws.Range("K12").Value = 1349
$ws.Range("L12").Value = 19309
$ws.Range("M12").Value = 5729
$ws.Range("N12").Value = 12319
$ws.Range("O12").Value = 18724
$ws.Range("P12").Value = 30587

# Problem 10
$ws.Range("B13").Value = 2649
$ws.Range("C13").Value = 401
$ws.Range("D13").Value = 2209
$ws.Range("E13").Value = 2496
$ws.Range("F13").Value = 3329
$ws.Range("G13").Value = 1063
$ws.Range("H13").Value = 98
$ws.Range("I13").Value = 985
$ws.Range("J13").Value = 1027
$ws.Range("K13").Value = 1262
$ws.Range("L13").Value = 16578
$ws.Range("M13").Value = 3136
$ws.Range("N13").Value = 12212
$ws.Range("O13").Value = 16433
$ws.Range("P13").Value = 21493

